# Add 2022-Q1 fund holdings sheet (copy of 2021-Q4 layout/styles), and
# update the summary ("总计") sheet with the new quarter's totals.
#
# Sheet order/ids in the target workbook are 2021-Q1(1), 2021-Q4(2),
# 2022-Q1(3), 总计(4) -- i.e. 总计 changes from sheetId 3 to 4 because a
# brand-new sheet (2022-Q1) is inserted before it. To reproduce that
# numbering we recreate 总计 (as a copy of 2021-Q4, reshaped into the
# totals layout) after the 2022-Q1 sheet has already taken the freed id.

$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $val) {
    # Force the cell to hold a literal text value (no auto-conversion to a
    # number) while leaving it with the "no explicit style" look of the
    # surrounding data cells.
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 0. Remove the existing "总计" sheet; we'll rebuild it at the end so
#    that sheet creation order/ids line up with the target workbook.
# ---------------------------------------------------------------------
$zjOld = $wb.Worksheets.Item("总计")
[void]$zjOld.Delete()

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet (holdings detail), positioned right after
#    "2021-Q4".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# The copied sheet has 3 data rows (like 2021-Q4); we only need 2, so
# drop the extra one before filling in the new values.
$q1.Rows.Item(4).Delete()

Set-TextCell $q1 2 2 "513090"
Set-TextCell $q1 2 3 "易方达中证香港证券投资主题ETF"
Set-TextCell $q1 2 4 "11.07"
Set-TextCell $q1 2 5 "96.47"
Set-TextCell $q1 2 6 "8.86"
Set-TextCell $q1 2 7 "0.9808"
$q1.Cells.Item(2, 8).Value = 4

Set-TextCell $q1 3 2 "002860"
Set-TextCell $q1 3 3 "前海开源沪港深新机遇灵活配置混合"
Set-TextCell $q1 3 4 "0.01"
Set-TextCell $q1 3 5 "83.26"
Set-TextCell $q1 3 6 "6.89"
Set-TextCell $q1 3 7 "0.0007"
$q1.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------
# 2. Rebuild "总计" (totals) sheet at the end, from another copy of
#    2021-Q4 (so it picks up the same header/index-column styling),
#    trimmed down to the 4-column totals layout.
# ---------------------------------------------------------------------
$q4b = $wb.Worksheets.Item("2021-Q4")
$q4b.Copy($null, $q1)
$zj = $wb.Worksheets.Item("2021-Q4 (2)")
$zj.Name = "总计"

# Drop the fund-specific columns E:H, keeping only A:D.
$zj.Range("E1:H4").EntireColumn.Delete()

# Header text is non-numeric, so a plain assignment keeps it as text
# without disturbing the inherited header style (s=2).
$zj.Cells.Item(1, 2).Value = "日期"
$zj.Cells.Item(1, 3).Value = "持有数量(只)"
$zj.Cells.Item(1, 4).Value = "持有市值(亿元)"

$zj.Cells.Item(2, 1).Value = 0
Set-TextCell $zj 2 2 "2022-Q1"
$zj.Cells.Item(2, 3).Value = 2
$zj.Cells.Item(2, 3).ClearFormats()
$zj.Cells.Item(2, 4).Value = 0.98
$zj.Cells.Item(2, 4).ClearFormats()

$zj.Cells.Item(3, 1).Value = 1
Set-TextCell $zj 3 2 "2021-Q4"
$zj.Cells.Item(3, 3).Value = 3
$zj.Cells.Item(3, 3).ClearFormats()
$zj.Cells.Item(3, 4).Value = 0.8
$zj.Cells.Item(3, 4).ClearFormats()

$zj.Cells.Item(4, 1).Value = 2
Set-TextCell $zj 4 2 "2021-Q1"
$zj.Cells.Item(4, 3).Value = 10
$zj.Cells.Item(4, 3).ClearFormats()
$zj.Cells.Item(4, 4).Value = 10.75
$zj.Cells.Item(4, 4).ClearFormats()

Write-Host "2022-Q1 sheet added and summary sheet rebuilt"
